# Add three more columns of sphere-diameter viscosity measurements
# (10 misurazioni con sferette di diametri diversi, distanza di 20 cm)
# to the "misure freddo" sheet, matching the layout already used on the
# "misure caldo" sheet (header row with 2,3,4,5,6 then numeric data
# formatted with two decimals).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("misure freddo")

# Header row: columns C, D, E get 4, 5, 6 (A=2, B=3 already present)
$ws.Range("C1").Value = 4
$ws.Range("D1").Value = 5
$ws.Range("E1").Value = 6

# Data values for rows 3-12, columns A through E
$data = @(
    @(21.79, 10.13, 6,    4.1,  2.96),
    @(21.5,  10.3,  6.02, 4.08, 2.96),
    @(21.75, 10.17, 6.07, 4.04, 2.96),
    @(22.04, 10.24, 5.91, 3.98, 2.95),
    @(21.71, 10.19, 6.07, 4.02, 2.95),
    @(21.77, 10.03, 6.03, 4.07, 3.13),
    @(22.03, 9.81,  5.84, 3.98, 2.91),
    @(21.79, 9.88,  5.84, 4.02, 2.83),
    @(21.44, 10.03, 6.04, 4.01, 2.89),
    @(21.66, 9.94,  5.84, 3.96, 2.95)
)

$row = 3
foreach ($values in $data) {
    $col = 1
    foreach ($v in $values) {
        $ws.Cells.Item($row, $col).Value = $v
        $col = $col + 1
    }
    $row = $row + 1
}

# Apply the same numeric formatting (two decimals) used on the analogous
# "misure caldo" sheet for this kind of measurement data.
$ws.Range("A3:E12").NumberFormat = "0.00"
